$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '67.837.47'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.36%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.621.24'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '595.09'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -1.61%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '152.73'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('E7').Value = '  +0.03%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.545'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.63%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.619.81'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('E10').Value = '  +6.72%  '
$ws.Range('E11').Value = '  -0.69%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '5.20'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('E13').Value = '  -1.19%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '27.55'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -1.76%  '
$ws.Range('E15').Value = '  +2.08%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.097.27'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -0.52%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '67.634.97'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.00%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.635.13'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.16%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '372.82'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.86%  '
$ws.Range('E20').Value = '  -0.59%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '7.43'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -2.61%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.23'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('E23').Value = '  -3.29%  '
$ws.Range('E24').Value = '  -4.80%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '72.55'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +9.50%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.87'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '592.57'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '0.0000103'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('E31').Value = '  -0.58%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '7.80'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -1.74%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.38'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -3.29%  '
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('E36').Value = '  -3.60%  '
$ws.Range('E37').Value = '  -1.15%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '158.28'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.14%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '19.12'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -1.91%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '1.90'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +3.12%  '
$ws.Range('E41').Value = '  -1.38%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.26'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('E43').Value = '  +2.13%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '17.12'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +4.56%  '
$ws.Range('E45').Value = '  +0.05%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '40.43'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -1.74%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '156.31'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.08%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.0₆0297'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +2.41%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '3.68'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('E50').Value = '  -2.59%  '
$ws.Range('E51').Value = '  -1.37%  '
